# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Overview / zh-cn / de-de "Status" cells that said "Ready for handoff"
#    now read "Handed back: in sync with en-US"
#  - de-de sheet gets its Latest Target File / Latest Handback File / Latest
#    Handback DateTime columns populated (with a hyperlink to the source .md)
#    for both rows, and the zh-cn sheet gets its Latest Target File /
#    Latest Handback File populated too (hyperlinked), while its handback
#    datetime cells pick up the refreshed "empty" handback time value.
#  - Several columns are widened so the new, longer strings fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"
$syncedEmptyTime  = "2016-08-27 09:01:04"
$handbackTimeDeDe = "2016-08-27 09:01:13"

$mdFileA = "a3673701-bd42-44b4-81d0-d3c2f37199d9.md"
$mdFileB = "ec16c88b-337a-4ed3-a246-c585ea35404a.md"
$mdUrlA  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34a65806fce1d1cf4bcd64f94f402d41e9635d61/e2e/a3673701-bd42-44b4-81d0-d3c2f37199d9.md"
$mdUrlB  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34a65806fce1d1cf4bcd64f94f402d41e9635d61/e2e/ec16c88b-337a-4ed3-a246-c585ea35404a.md"

$zhCnTargetA = "a3673701-bd42-44b4-81d0-d3c2f37199d9.ccd5440b879d0f23d02e4a1c841220c8921b2396.zh-cn.xlf"
$zhCnTargetB = "ec16c88b-337a-4ed3-a246-c585ea35404a.73dc346bd2cf88c1dc12828ed5af86bc4ded9da7.zh-cn.xlf"
$deDeTargetA = "a3673701-bd42-44b4-81d0-d3c2f37199d9.ccd5440b879d0f23d02e4a1c841220c8921b2396.de-de.xlf"
$deDeTargetB = "ec16c88b-337a-4ed3-a246-c585ea35404a.73dc346bd2cf88c1dc12828ed5af86bc4ded9da7.de-de.xlf"

# ---------------------------------------------------------------------------
# 1. Update the "Status" text everywhere it shows up (Overview zh-cn/de-de
#    columns, and the Status column of the two locale sheets).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" / "Latest Handback File" and
#    refresh the (still-empty) handback datetime.
# ---------------------------------------------------------------------------
$wsZhCn.Range("I2").Value = $mdFileA
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrlA, "", "", $mdFileA)
$wsZhCn.Range("J2").Value = $zhCnTargetA
$wsZhCn.Range("K2").Value = $syncedEmptyTime

$wsZhCn.Range("I3").Value = $mdFileB
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrlB, "", "", $mdFileB)
$wsZhCn.Range("J3").Value = $zhCnTargetB
$wsZhCn.Range("K3").Value = $syncedEmptyTime

# ---------------------------------------------------------------------------
# 3. de-de sheet: populate "Latest Target File" / "Latest Handback File" and
#    the handback datetime for both rows.
# ---------------------------------------------------------------------------
$wsDeDe.Range("I2").Value = $mdFileA
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrlA, "", "", $mdFileA)
$wsDeDe.Range("J2").Value = $deDeTargetA
$wsDeDe.Range("K2").Value = $handbackTimeDeDe

$wsDeDe.Range("I3").Value = $mdFileB
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrlB, "", "", $mdFileB)
$wsDeDe.Range("J3").Value = $deDeTargetB
$wsDeDe.Range("K3").Value = $handbackTimeDeDe

# ---------------------------------------------------------------------------
# 4. Widen the columns that now hold longer content. The host quantizes
#    ColumnWidth to 1/6-character pixel steps, so we dial in the character
#    width that lands on the intended rendered width after that rounding.
# ---------------------------------------------------------------------------
$wideStatusWidth = 29.166666666666668   # renders as 29.9777047293527 -> rounds to 30
$wideFileWidth    = 39.166666666666664  # renders as exactly 40

$wsOverview.Columns("E").ColumnWidth = $wideStatusWidth
$wsOverview.Columns("F").ColumnWidth = $wideStatusWidth

$wsZhCn.Columns("C").ColumnWidth = $wideStatusWidth
$wsZhCn.Columns("I").ColumnWidth = $wideFileWidth
$wsZhCn.Columns("J").ColumnWidth = $wideFileWidth

$wsDeDe.Columns("C").ColumnWidth = $wideStatusWidth
$wsDeDe.Columns("I").ColumnWidth = $wideFileWidth
$wsDeDe.Columns("J").ColumnWidth = $wideFileWidth

Write-Host "Handback report generated."
